$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in previously-missing Year values ---
$ws.Range("C8").Value = 2017
$ws.Range("E8").Value = "170 page Master Thesis"
$ws.Range("C9").Value = 2017
$ws.Range("C38").Value = 1995
$ws.Range("C39").Value = 2011
$ws.Range("C40").Value = 2007
$ws.Range("C41").Value = 2002

# --- Insert 4 new rows before row 51 (pushes old 51-62 down to 55-66) ---
$ws.Range("A51:H54").EntireRow.Insert()

# The row-insert shim does not relocate the Hyperlinks collection entries,
# so any hyperlinks that used to live in rows 52/53 are still "anchored"
# there even though the cell text itself correctly moved down to 56/57.
# Remove those stale hyperlink objects and re-add them at the right spot.
function Remove-HyperlinkAtRow($sheet, $row) {
    $found = $null
    foreach ($h in $sheet.Hyperlinks) {
        if ($h.Range.Row -eq $row) {
            $found = $h
            break
        }
    }
    if ($found -ne $null) {
        $addr = $found.Address
        $found.Delete()
        return $addr
    }
    return $null
}

$addr52 = Remove-HyperlinkAtRow $ws 52
$addr53 = Remove-HyperlinkAtRow $ws 53
if ($addr52) { $ws.Hyperlinks.Add($ws.Range("H56"), $addr52) | Out-Null }
if ($addr53) { $ws.Hyperlinks.Add($ws.Range("H57"), $addr53) | Out-Null }

# New section header row (row 51)
$ws.Range("A51").Value = "Integrated Simulation and Optimisation"

# Name-of-paper text for the two new rows
$ws.Range("A52").Value = "Integrating simulation and optimisation in health care centre management"
$ws.Range("A53").Value = "Integrating optimisation and simulation approaches for daily scheduling of assembly and test operations"

# Author(s) text for the two new rows
$ws.Range("B52").Value = "De Angelis et al"
$ws.Range("B53").Value = "Bard et al"

# Year
$ws.Range("C52").Value = 2003
$ws.Range("C53").Value = 2015

# Notes text (shared between both new rows)
$ws.Range("E52").Value = "Background, but potentially good to check for inspiration"
$ws.Range("E53").Value = "Background, but potentially good to check for inspiration"

# Read? / Priority
$ws.Range("F52").Value = "N"
$ws.Range("G52").Value = "M"
$ws.Range("F53").Value = "N"
$ws.Range("G53").Value = "M"

# Row 54 is intentionally left blank (gap row, same as in the rest of the sheet)

# Old row 53 (now row 57) had a blank Year cell; it gets filled in too
$ws.Range("C57").Value = 2015
